# Update FFXIV Leve profit-calculation columns (H-N) across several sheets
# per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 325.8095
$ws.Range("I38").Value = 200.125
$ws.Range("J38").Value = 728
$ws.Range("K38").Value = 600.375
$ws.Range("L38").Value = 2184
$ws.Range("M38").Value = -228.375
$ws.Range("N38").Value = -2928

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1756.75
$ws.Range("I43").Value = 925
$ws.Range("J43").Value = 2034
$ws.Range("K43").Value = 925
$ws.Range("L43").Value = 2034
$ws.Range("M43").Value = -856
$ws.Range("N43").Value = -2172

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 948.5833
$ws.Range("I58").Value = 948.5833
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2845.7499
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2695.7499
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 998.2857
$ws.Range("I125").Value = 1131.6666
$ws.Range("J125").Value = 944.93335
$ws.Range("K125").Value = 10184.9994
$ws.Range("L125").Value = 8504.400149999999
$ws.Range("M125").Value = -7724.999400000001
$ws.Range("N125").Value = -13424.40015

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2913.0222
$ws.Range("I137").Value = 2002.6216
$ws.Range("J137").Value = 7123.625
$ws.Range("K137").Value = 6007.864799999999
$ws.Range("L137").Value = 21370.875
$ws.Range("M137").Value = -3457.864799999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1869.38
$ws.Range("I138").Value = 586.75
$ws.Range("J138").Value = 2113.6904
$ws.Range("K138").Value = 1760.25
$ws.Range("L138").Value = 6341.0712
$ws.Range("M138").Value = 3379.75
$ws.Range("N138").Value = -16621.0712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3267.348
$ws.Range("I61").Value = 2704.5
$ws.Range("J61").Value = 3567.5334
$ws.Range("K61").Value = 2704.5
$ws.Range("L61").Value = 3567.5334
$ws.Range("M61").Value = -2492.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2538.6216
$ws.Range("I74").Value = 2309.2104
$ws.Range("J74").Value = 2780.7778
$ws.Range("K74").Value = 2309.2104
$ws.Range("L74").Value = 2780.7778
$ws.Range("M74").Value = -1435.2104
$ws.Range("N74").Value = -4528.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2538.6216
$ws.Range("I77").Value = 2309.2104
$ws.Range("J77").Value = 2780.7778
$ws.Range("K77").Value = 11546.052
$ws.Range("L77").Value = 13903.889
$ws.Range("M77").Value = -7178.052
$ws.Range("N77").Value = -22639.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3689.1738
$ws.Range("I132").Value = 4231.048
$ws.Range("J132").Value = 3234
$ws.Range("K132").Value = 12693.144
$ws.Range("L132").Value = 9702
$ws.Range("M132").Value = -10163.144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3267.348
$ws.Range("I136").Value = 2704.5
$ws.Range("J136").Value = 3567.5334
$ws.Range("K136").Value = 8113.5
$ws.Range("L136").Value = 10702.6002
$ws.Range("M136").Value = -5563.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2324.8572
$ws.Range("I134").Value = 2083.8667
$ws.Range("J134").Value = 2927.3333
$ws.Range("K134").Value = 6251.6001
$ws.Range("L134").Value = 8781.999899999999
$ws.Range("M134").Value = -3716.6001
$ws.Range("N134").Value = -13851.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 10500
$ws.Range("I26").Value = 1000
$ws.Range("J26").Value = 20000
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = -713
$ws.Range("N26").Value = -20574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 204135
$ws.Range("I35").Value = 204135
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 204135
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -203841

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1756.9474
$ws.Range("I58").Value = 1732
$ws.Range("J58").Value = 1771.5
$ws.Range("K58").Value = 1732
$ws.Range("L58").Value = 1771.5
$ws.Range("M58").Value = -1529
$ws.Range("N58").Value = -2177.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 27780196
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 41669292
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 125007876
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -125012936

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2535.84
$ws.Range("I134").Value = 1672.8334
$ws.Range("J134").Value = 4755
$ws.Range("K134").Value = 5018.5002
$ws.Range("L134").Value = 14265
$ws.Range("M134").Value = -2483.5002
$ws.Range("N134").Value = -19335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1756.9474
$ws.Range("I136").Value = 1732
$ws.Range("J136").Value = 1771.5
$ws.Range("K136").Value = 5196
$ws.Range("L136").Value = 5314.5
$ws.Range("M136").Value = -2646
$ws.Range("N136").Value = -10414.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1584.5071
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1584.5071
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 4753.5213
$ws.Range("N39").Value = -5341.5213

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 11534.333
$ws.Range("I110").Value = 4342.3335
$ws.Range("J110").Value = 12253.533
$ws.Range("K110").Value = 13027.0005
$ws.Range("L110").Value = 36760.599
$ws.Range("M110").Value = -8937.000499999998
$ws.Range("N110").Value = -44940.599

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 997
$ws.Range("I131").Value = 460
$ws.Range("J131").Value = 1045.8182
$ws.Range("K131").Value = 1380
$ws.Range("L131").Value = 3137.4546
$ws.Range("M131").Value = 3660
$ws.Range("N131").Value = -13217.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2675.963
$ws.Range("I139").Value = 2031.875
$ws.Range("J139").Value = 2947.158
$ws.Range("K139").Value = 6095.625
$ws.Range("L139").Value = 8841.474
$ws.Range("M139").Value = -955.625
$ws.Range("N139").Value = -19121.474

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3325.889
$ws.Range("I132").Value = 3030.2856
$ws.Range("J132").Value = 3514
$ws.Range("K132").Value = 9090.856800000001
$ws.Range("L132").Value = 10542
$ws.Range("M132").Value = -6560.856800000001
$ws.Range("N132").Value = -15602

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7978.6665
$ws.Range("I22").Value = 748.1667
$ws.Range("J22").Value = 10044.523
$ws.Range("K22").Value = 748.1667
$ws.Range("L22").Value = 10044.523
$ws.Range("M22").Value = -453.1667
$ws.Range("N22").Value = -10634.523

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 7978.6665
$ws.Range("I27").Value = 748.1667
$ws.Range("J27").Value = 10044.523
$ws.Range("K27").Value = 748.1667
$ws.Range("L27").Value = 10044.523
$ws.Range("M27").Value = -641.1667
$ws.Range("N27").Value = -10258.523

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5040
$ws.Range("I46").Value = 1233.3334
$ws.Range("J46").Value = 10750
$ws.Range("K46").Value = 1233.3334
$ws.Range("L46").Value = 10750
$ws.Range("M46").Value = -1045.3334
$ws.Range("N46").Value = -11126

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 678.5
$ws.Range("I55").Value = 475.25
$ws.Range("J55").Value = 729.3125
$ws.Range("K55").Value = 475.25
$ws.Range("L55").Value = 729.3125
$ws.Range("M55").Value = -302.25
$ws.Range("N55").Value = -1075.3125

Write-Output "Updated 27 rows (182 cell writes, 1 clears)"
